# Applies the cryptos list update described in the commit diff.
# Cell values are plain text (inline strings) in the source workbook;
# for values that look like plain numbers we prefix with an apostrophe
# so Excel stores them as text (matching the original inlineStr cells)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.414.55"
$ws.Range("D3").Value = "1.801.68"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'227.71"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'0.581"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'35.00"
$ws.Range("E8").Value = "  +5.99%  "
$ws.Range("D9").Value = "'0.299"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "2.060.65"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'11.18"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "1.797.36"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'0.641"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "34.363.83"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'244.74"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'4.15"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'170.81"
$ws.Range("E24").Value = "  +3.50%  "
$ws.Range("D25").Value = "'2.11"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "'7.56"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("D27").Value = "'16.72"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0527"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "1.399.77"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "'82.94"
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("D42").Value = "'0.947"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'13.64"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "'5.98"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "1.962.45"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "'104.51"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  +0.08%  "
